$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "-"
$ws.Range("B3").Value = "MCT-3A-Máquinas Elétricas"

$ws.Range("B4").Value = "MCT-3A-Máquinas Elétricas"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "MCT-1A-Circuitos elétricos"

$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "MCT-1A-Circuitos elétricos"

$ws.Range("C7").Value = "-"
$ws.Range("E7").Value = "MCT-1A-Circuitos elétricos"

$ws.Range("B11").Value = "-"
$ws.Range("B15").Value = "Cleidson-Circuitos elétricos-1A"

$ws.Range("D18").Value = "-"

$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "['ELM-1NA-Sistemas digitais', 'ELM-2NA-Automação Industrial', -, -]"
$ws.Range("F19").Value = "[-, -, 'ELM-2NA-Automação Industrial', -]"

$ws.Range("B20").Value = "-"
$ws.Range("F20").Value = "[-, -, 'ELM-2NA-Automação Industrial', -]"

$ws.Range("B21").Value = "-"
$ws.Range("F21").Value = "[-, -, 'ELM-2NA-Automação Industrial', -]"
